# Remove the three paragraphs that followed the "LOQ4084: Fenomenos de
# Transporte II (Requisito fraco)" requirement line:
#   - an empty paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#     pages. Original theme under Creative Commons Attribution"
# leaving the requirement paragraph directly followed by the (previously
# trailing) empty paragraph right before the page break.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the copyright paragraph by
# their text, so the script does not depend on fixed paragraph indices.
$count = $d.Paragraphs.Count
$verIndex = -1
$copyrightIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Ver no Jupiter*") {
        $verIndex = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

# The blank paragraph right before "Ver no Jupiter ..." is also removed.
$blankIndex = $verIndex - 1

$rangeStart = $d.Paragraphs.Item($blankIndex).Range.Start
$rangeEnd = $d.Paragraphs.Item($copyrightIndex).Range.End

$r = $d.Range($rangeStart, $rangeEnd)
$r.Delete()
